$x = 1 + 2
Write-Host "sum:" $x
Write-Host "test env"
Write-Host ([System.IO.Path]::GetTempPath())
